$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(18, 19),
    @(52, 53),
    @(105, 106),
    @(114, 115),
    @(130, 131),
    @(133, 134),
    @(150, 151),
    @(158, 160),
    @(197, 198),
    @(203, 204),
    @(210, 211),
    @(226, 227),
    @(229, 230)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
